# The deck's Design/theme was switched from the custom "Integral" colour
# scheme to the stock "Office Theme" colour scheme (PowerPoint Design tab
# -> Office Theme). The font scheme (Arial everywhere) and the effect/fill
# format scheme are already identical between the two themes, so the only
# substantive change is the 12 theme colours (and the theme's display
# name). Re-apply the standard "Office" theme palette to the presentation's
# slide-master theme so every placeholder/shape that references a scheme
# colour (bg1/tx1/bg2/tx2/accent1-6/hlink/folHlink) picks up the new
# colours, matching ppt/theme/theme1.xml in the target deck.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme
$cs = $theme.ThemeColorScheme

# Index -> (name, RGB) using the MsoThemeColorSchemeIndex ordering:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$cs.Item(1).RGB  = 0        # dk1      000000
$cs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388  # dk2      44546A
$cs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501  # accent2  ED7D31
$cs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$cs.Item(8).RGB  = 49407    # accent4  FFC000
$cs.Item(9).RGB  = 12874308 # accent5  4472C4
$cs.Item(10).RGB = 4697456  # accent6  70AD47
$cs.Item(11).RGB = 12673797 # hlink    0563C1
$cs.Item(12).RGB = 7491477  # folHlink 954F72

# Best-effort: also rename the theme/colour-scheme to "Office"/"Office
# Theme" (some hosts expose these as read-only computed names, so guard
# each assignment individually and keep going either way).
try { $theme.Name = "Office Theme" } catch { }
try { $cs.Name = "Office" } catch { }
try { $p.Designs.Item(1).Name = "Office Theme" } catch { }
